$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260903716087341
$ws.Range("B1").Value = 1.512049674987793
$ws.Range("C1").Value = 2.017351865768433
$ws.Range("D1").Value = 2.106147527694702
$ws.Range("E1").Value = 1.166621208190918
